$p = $ppt.ActivePresentation

# 1. Update the cached "datetimeFigureOut" field text from 7/20/17 to 4/16/2018
#    across the slide master, all slide layouts, and the notes master.
function Set-DateText($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "7/20/17") {
                $sh.TextFrame.TextRange.Text = "4/16/2018"
            }
        }
    }
}

Set-DateText $p.SlideMaster
for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    Set-DateText $p.SlideMaster.CustomLayouts.Item($j)
}
Set-DateText $p.NotesMaster

# 2. Remove the now-obsolete "UndoRedo Stack" diagram (rectangle, connector arrow,
#    and its "1" textbox label) from slide 1 - this mini-diagram is no longer
#    relevant now that LogicManager no longer owns the undo/redo stack.
$s = $p.Slides.Item(1)
$s.Shapes.Item(51).Delete()   # TextBox 62 ("1")
$s.Shapes.Item(50).Delete()   # Straight Arrow Connector 57
$s.Shapes.Item(49).Delete()   # Rectangle 62 ("UndoRedo Stack")
